# Increment the "Förändrad" (Changed) date column C by one day for every
# data row (rows 2 through 250), changing the serial value from 45177 to 45178.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 250
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
